$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (for line7 and line8) after the existing line6 row,
# pushing the extr1..extr8 rows down by two rows.
$ws.Range("A8:A9").EntireRow.Insert()

# Match the formatting used by the other index cells in column A
# (bold font, thin border on all sides, centered horizontal/top vertical alignment).
$ws.Range("A8:A9").Font.Bold = $true
$ws.Range("A8:A9").Borders.LineStyle = 1
$ws.Range("A8:A9").HorizontalAlignment = -4108
$ws.Range("A8:A9").VerticalAlignment = -4160

# New row: line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row: line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $false

# Renumber the index column for the extr1..extr8 rows, which have been
# pushed down to rows 10..17, and flip in_service for extr1 and extr2.
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

$ws.Range("E10").Value = $true
$ws.Range("E11").Value = $true
